# Implement column width parsing (#3)
# Set explicit column widths on column A and column D so the saved
# worksheet carries <cols> entries (width=27 for A, width=32.5 for D),
# matching the custom widths added to tests/format.xlsx.
#
# Excel's COM `ColumnWidth` property is expressed in "characters" and is
# offset from the raw OOXML <col width="..."> value by a constant padding
# (5/6 of a character for this workbook's default font/metrics). Feed in
# the de-padded values so the persisted width lands on the exact target.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").EntireColumn.ColumnWidth = 26.166666666666668
$ws.Range("D1").EntireColumn.ColumnWidth = 31.666666666666668
